$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the restated values for the existing last three rows (314-316) ---
$ws.Cells.Item(314, 3).Value = 2078216858000
$ws.Cells.Item(314, 4).Value = 2078216858000
$ws.Cells.Item(314, 5).Value = 2078216858000
$ws.Cells.Item(314, 6).Value = 2078216858000

$ws.Cells.Item(315, 3).Value = 2082183969000
$ws.Cells.Item(315, 4).Value = 2082183969000
$ws.Cells.Item(315, 5).Value = 2082183969000
$ws.Cells.Item(315, 6).Value = 2082183969000

$ws.Cells.Item(316, 3).Value = 2118202312000
$ws.Cells.Item(316, 4).Value = 2118202312000
$ws.Cells.Item(316, 5).Value = 2118202312000
$ws.Cells.Item(316, 6).Value = 2118202312000

# --- Append three new monthly data rows (317-319) ---
$newRows = @(
    @{ Row = 317; Date = 44986.45833333334; Value = 2121975670000 },
    @{ Row = 318; Date = 45017.45833333334; Value = 2135028350000 },
    @{ Row = 319; Date = 45047.41666666666; Value = 2140971740000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the date-cell formatting (style index carrying the datetime
    # number format/border/font) from the row above down to the new row.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:PLM2"
    $ws.Cells.Item($row, 3).Value = $r.Value
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 5).Value = $r.Value
    $ws.Cells.Item($row, 6).Value = $r.Value
    $ws.Cells.Item($row, 7).Value = 0
}

$excel.CutCopyMode = $false
